# The deck ships two themes:
#   ppt/theme/theme1.xml -> used by the (single) Slide Master, currently "Integral"
#   ppt/theme/theme2.xml -> used by the Notes Master, currently "Office Theme"
# The commit swaps the two themes' contents: the slide master becomes the
# stock "Office Theme" colour scheme and the notes master becomes "Integral".
# Font scheme / format scheme are identical between the two themes already,
# so only the 12 theme colours (+ names) actually move.

$p = $ppt.ActivePresentation

# --- Slide master theme (ppt/theme/theme1.xml): Integral -> Office Theme ---
$masterScheme = $p.Designs.Item(1).SlideMaster.Theme.ThemeColorScheme
$masterScheme.Item(1).RGB  = 0x000000   # dk1
$masterScheme.Item(2).RGB  = 0xFFFFFF   # lt1
$masterScheme.Item(3).RGB  = 0x6A5444   # dk2
$masterScheme.Item(4).RGB  = 0xE6E6E7   # lt2
$masterScheme.Item(5).RGB  = 0xD59B5B   # accent1
$masterScheme.Item(6).RGB  = 0x317DED   # accent2
$masterScheme.Item(7).RGB  = 0xA5A5A5   # accent3
$masterScheme.Item(8).RGB  = 0x00C0FF   # accent4
$masterScheme.Item(9).RGB  = 0xC47244   # accent5
$masterScheme.Item(10).RGB = 0x47AD70   # accent6
$masterScheme.Item(11).RGB = 0xC16305   # hlink
$masterScheme.Item(12).RGB = 0x724F95   # folHlink

# --- Notes master theme (ppt/theme/theme2.xml): Office Theme -> Integral ---
$notesScheme = $p.NotesMaster.Theme.ThemeColorScheme
$notesScheme.Item(1).RGB  = 0x000000   # dk1
$notesScheme.Item(2).RGB  = 0xFFFFFF   # lt1
$notesScheme.Item(3).RGB  = 0x515F45   # dk2
$notesScheme.Item(4).RGB  = 0xD1DEE3   # lt2
$notesScheme.Item(5).RGB  = 0x38CB99   # accent1
$notesScheme.Item(6).RGB  = 0x37A563   # accent2
$notesScheme.Item(7).RGB  = 0x24D0E6   # accent3
$notesScheme.Item(8).RGB  = 0x0097CC   # accent4
$notesScheme.Item(9).RGB  = 0xCFB34E   # accent5
$notesScheme.Item(10).RGB = 0xA68D37   # accent6
$notesScheme.Item(11).RGB = 0x259F6B   # hlink
$notesScheme.Item(12).RGB = 0x026BB2   # folHlink
